{"js": "// Office.js (Word JavaScript API) script implementing:\n//  1. Insert a new \"Meta description\" paragraph right after the document's\n//     title (Heading1) paragraph.\n//  2. Remove the duplicate bolded title paragraph near the end of the\n//     document (it is superseded by the new meta-description block).\n//  3. Replace the text of the trailing italic paragraph (previously the\n//     meta description copy) with the new image-generation prompt text,\n//     keeping its italic formatting intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Step 1: insert the \"Meta description\" paragraph after the title ----\nconst titlePara = paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nmetaPara.style = \"Normal\";\nawait context.sync();\n\nconst boldText = \"Meta description\";\nconst restText =\n  \": Read our review of Cyberslot Megaclusters, an innovative online slot game with unique gameplay mechanics and multiplier features. Play for free here.\";\n\n// Insert the non-bold text first (at the paragraph's end), then insert the\n// bold lead-in before it - this keeps the two segments as separate runs\n// without leaving a stray empty rPr on the second run.\nmetaPara.insertText(restText, \"End\");\nawait context.sync();\n\nconst boldRange = metaPara.insertText(boldText, \"Start\");\nboldRange.font.bold = true;\nawait context.sync();\n\n// --- Step 2: delete the duplicate bold title paragraph near the end -----\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\n\nconst items = allParagraphs.items;\nconst dupTitlePara = items[items.length - 2];\ndupTitlePara.delete();\nawait context.sync();\n\n// --- Step 3: update the trailing italic paragraph's text -----------------\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\n\nconst lastPara = remaining.items[remaining.items.length - 1];\nconst promptText =\n  \"Prompt: Create a colorful and eye-catching feature image in a cartoon style for Cyberslot Megaclusters. The image should prominently feature a happy Maya warrior wearing glasses, in line with the futuristic and technology-themed game. The image should convey the game's mini-grid mechanic and use a bright color palette to reflect the game's simple yet modern aesthetic. Please include the game title and any additional elements that you feel would enhance the image's appeal and accurately represent the game's features.\";\n\nlastPara.getRange(\"Whole\").insertText(promptText, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop script implementing:\n#  1. Insert a new \"Meta description\" paragraph right after the document's\n#     title (Heading1) paragraph.\n#  2. Remove the duplicate bolded title paragraph near the end of the\n#     document (it is superseded by the new meta-description block).\n#  3. Replace the text of the trailing italic paragraph (previously the\n#     meta description copy) with the new image-generation prompt text,\n#     keeping its italic formatting intact.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the \"Meta description\" paragraph after the title ----\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter() | Out-Null\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n\n$boldText = \"Meta description\"\n$restText = \": Read our review of Cyberslot Megaclusters, an innovative online slot game with unique gameplay mechanics and multiplier features. Play for free here.\"\n\n$metaRange = $metaPara.Range\n$metaRange.Collapse(1)\n$metaRange.InsertAfter($boldText + $restText)\n\n$metaStart = $metaPara.Range.Start\n$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)\n$boldRange.Bold = 1\n\n# --- Step 2: delete the duplicate bold title paragraph near the end -----\n$count = $d.Paragraphs.Count\n$dupTitlePara = $d.Paragraphs.Item($count - 1)\n$dupTitlePara.Range.Delete()\n\n# --- Step 3: update the trailing italic paragraph's text -----------------\n$count = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($count)\n$promptText = \"Prompt: Create a colorful and eye-catching feature image in a cartoon style for Cyberslot Megaclusters. The image should prominently feature a happy Maya warrior wearing glasses, in line with the futuristic and technology-themed game. The image should convey the game's mini-grid mechanic and use a bright color palette to reflect the game's simple yet modern aesthetic. Please include the game title and any additional elements that you feel would enhance the image's appeal and accurately represent the game's features.\"\n\n$lastRange = $lastPara.Range\n$lastRange.End = $lastRange.End - 1\n$lastRange.Text = $promptText\n"}
